$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header row: question labels
$ws.Range("B1").Value = "Question 1"
$ws.Range("C1").Value = "Question 2"

# New second row: quiz label + actual question text (answers removed)
$ws.Range("A2").Value = "Excel-quiz"
$ws.Range("B2").Value = "Where do you find the best answers?"
$ws.Range("C2").Value = "Who to ask?"
